$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("What is House Renting?") - Content Placeholder 2
#   * "a process that" -> "a system that"
#   * drop the trailing "It is reduce time and cost." bullet (moved to
#     slide 9 with new wording later in the deck)
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Text = "House renting is a system that is People can advertising his or her own house or People can view house advertisement who want to rent house." + [char]13 + "People can easily find expected house for rent."

# ---------------------------------------------------------------------------
# Slide 5 ("Feasibility Study") - Content Placeholder 2
#   * shorten the intro sentence
#   * add detail to "Technical Feasibility: "
#   * add a trailing space to "Economic Feasibility:"
#   * merge the two "Environmental" / ":" runs into a single run
# paragraph count is unchanged (4), so edit via Characters() ranges that
# exactly span each whole paragraph/run to keep per-paragraph pPr (lvl/sz).
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange

# Paragraph 4 has two runs: "Environmental" (chars 157-169) + ":" (char 170).
# Write the merged text into the *second* run first (so the merged run
# keeps that run's rPr, which carries dirty="0"), then clear the first run.
$p4run2 = $tr5.Characters(170, 1)
$p4run2.Text = "Environmental:"
$p4run1 = $tr5.Characters(157, 13)
$p4run1.Text = ""

# Paragraph 3 (single run): "Economic Feasibility:" (chars 135-155)
$p3 = $tr5.Characters(135, 21)
$p3.Text = "Economic Feasibility: "

# Paragraph 2 (single run): "Technical Feasibility: " (chars 111-133)
$p2 = $tr5.Characters(111, 23)
$p2.Text = "Technical Feasibility: I have 1year experience on angular for font end design and 6month experience on node.js framework for back end design."

# Paragraph 1 (single run): chars 1-109
$p1 = $tr5.Characters(1, 109)
$p1.Text = "A feasibility study evaluates the project" + [char]0x2019 + "s potential for success."

# ---------------------------------------------------------------------------
# Slide 6 ("Application Area") - Content Placeholder 2
#   * "Academic:" -> "Actually this project is popular in urban area."
#   * add a new bullet "All class of professionals will be most popularities."
# paragraph count changes (2 -> 3), so rewrite the whole text range.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange
$tr6.Text = "Industrial: This Project is used as business purposes. " + [char]13 + "Actually this project is popular in urban area." + [char]13 + "All class of professionals will be most popularities."

# ---------------------------------------------------------------------------
# Slide 9 ("Expected Outcome") - Content Placeholder 2
#   * fill in the empty trailing paragraph with real text
#   * add one more new bullet paragraph
# paragraph count changes (2 -> 3), so rewrite the whole text range.
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange
$tr9.Text = "Complete business Platform." + [char]13 + "It is reduce time and cost for both renter and advertiser." + [char]13 + "It is help to hurry up to find  a house for rent."

# ---------------------------------------------------------------------------
# Slide 10 ("Risk of Our Project") - Content Placeholder 2
#   * "to take information, subscription" -> "to take information and subscription"
# paragraph count unchanged; replace the whole (single-run) paragraph text.
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(2).TextFrame.TextRange
$p10last = $tr10.Characters(179, 86)
$p10last.Text = "Sometimes bad people to take information and subscription is apply on house advertiser."
